# Horarios actualizados Línea 141 - 954
# Updates the scrape timestamp / minute countdowns for the three schedule
# sheets (LP1912, LP1912-215, 6203-6173) and appends the newly-scraped
# rows that came in with this refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 06:56:24"
$ws.Range("A3").Value = "Total filas: 49"

# Rows 16 & 17 swap their Hora_Scrap / Linea / Minutos values.
$ws.Range("A16").Value = "06:38:54"
$ws.Range("C16").Value = "16_SANTA ANA"
$ws.Range("D16").Value = 2

$ws.Range("A17").Value = "05:44:02"
$ws.Range("C17").Value = "17X38_ROMERO"
$ws.Range("D17").Value = 56

$ws.Range("A20").Value = "06:56:24"
$ws.Range("D20").Value = 1

$ws.Range("A22").Value = "06:56:24"
$ws.Range("D22").Value = 3

$ws.Range("A24").Value = "06:56:24"
$ws.Range("D24").Value = 20

$ws.Range("A26").Value = "06:56:24"
$ws.Range("D26").Value = 23

$ws.Range("A28").Value = "06:56:24"
$ws.Range("C28").Value = "16_SANTA ANA"
$ws.Range("D28").Value = 25

# Row 29 swaps its Linea value with row 28 (A/B/D/E stay the same).
$ws.Range("C29").Value = "23_HERNANDEZ"

$ws.Range("A30").Value = "06:56:24"
$ws.Range("D30").Value = 26

$ws.Range("A31").Value = "06:56:24"
$ws.Range("D31").Value = 33

$ws.Range("A33").Value = "06:56:24"
$ws.Range("D33").Value = 39

$ws.Range("A35").Value = "06:56:24"
$ws.Range("D35").Value = 41

$ws.Range("A37").Value = "06:56:24"
$ws.Range("D37").Value = 48

$ws.Range("A39").Value = "06:56:24"
$ws.Range("D39").Value = 59

$ws.Range("A42").Value = "06:56:24"
$ws.Range("D42").Value = 64

$ws.Range("A43").Value = "06:56:24"
$ws.Range("D43").Value = 65

$ws.Range("A45").Value = "06:56:24"
$ws.Range("D45").Value = 75

$ws.Range("A47").Value = "06:56:24"
$ws.Range("D47").Value = 77

$ws.Range("A49").Value = "06:56:24"
$ws.Range("D49").Value = 93

$ws.Range("A50").Value = "06:56:24"
$ws.Range("D50").Value = 93

$ws.Range("A51").Value = "06:56:24"
$ws.Range("D51").Value = 105

$ws.Range("A52").Value = "06:56:24"
$ws.Range("D52").Value = 108

# Two brand-new rows appended at the bottom of the sheet.
$ws.Range("A53").Value = "06:56:24"
$ws.Range("B53").Value = "08:52"
$ws.Range("C53").Value = "23_HERNANDEZ"
$ws.Range("D53").Value = 116
$ws.Range("E53").Value = "LP1912"

$ws.Range("A54").Value = "06:56:24"
$ws.Range("B54").Value = "08:53"
$ws.Range("C54").Value = "215B_EL PATO"
$ws.Range("D54").Value = 117
$ws.Range("E54").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A2").Value = "Última actualización: 06:56:24"
$ws.Range("A3").Value = "Total filas: 9"

$ws.Range("A8").Value = "06:56:24"
$ws.Range("D8").Value = 1

$ws.Range("A10").Value = "06:56:24"
$ws.Range("D10").Value = 20

$ws.Range("A12").Value = "06:56:24"
$ws.Range("D12").Value = 48

$ws.Range("A13").Value = "06:56:24"
$ws.Range("D13").Value = 108

# New row appended at the bottom of the sheet.
$ws.Range("A14").Value = "06:56:24"
$ws.Range("B14").Value = "08:53"
$ws.Range("C14").Value = "215B_EL PATO"
$ws.Range("D14").Value = 117
$ws.Range("E14").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A2").Value = "Última actualización: 06:56:24"
$ws.Range("A3").Value = "Total filas: 5"

$ws.Range("A7").Value = "06:56:24"
$ws.Range("D7").Value = 47

$ws.Range("A9").Value = "06:56:24"
$ws.Range("D9").Value = 100

# New row appended at the bottom of the sheet.
$ws.Range("A10").Value = "06:56:24"
$ws.Range("B10").Value = "08:51"
$ws.Range("C10").Value = "215C_LA PLATA"
$ws.Range("D10").Value = 115
$ws.Range("E10").Value = "L6203"
